$d = $word.ActiveDocument

# ===========================================================================
# PHASE 1 - structural changes (paragraph insert/delete), done first and
# purely by paragraph index so nothing here depends on a "live" reference
# surviving a later structural edit.
# ===========================================================================

# 1a) Remove the "Outra etapa importante do projeto..." paragraph entirely
#     (paragraph 59 in the original document).
$d.Paragraphs.Item(59).Range.Delete()

# After the delete above the "Escopo" section paragraphs are:
#   55 Escopo (heading)
#   56 O meu projeto tem como objetivo...
#   57 A primeira etapa...
#   58 Além disso, o projeto inclui a criação de um quiz...
#   59 Além disso, este projeto também terá uma seção de biblioteca virtual...
#   60 Premissas e Restrições (heading)

# 1b) Insert a new empty paragraph right before paragraph 56 ("O meu
#     projeto..."). It will become paragraph 56 and push the rest down by one.
$r = $d.Paragraphs.Item(56).Range
$r.Collapse(1)
$r.InsertParagraphBefore()

# Now:
#   55 Escopo (heading)
#   56 <new blank paragraph>
#   57 O meu projeto tem como objetivo...
#   58 A primeira etapa...
#   59 Além disso, o projeto inclui a criação de um quiz...
#   60 Além disso, este projeto também terá uma seção de biblioteca virtual...
#   61 Premissas e Restrições (heading)

# 1c) Insert a new empty paragraph right after paragraph 60 (the "biblioteca
#     virtual" paragraph).
$r2 = $d.Paragraphs.Item(60).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# Now:
#   55 Escopo (heading)
#   56 <new blank paragraph>
#   57 O meu projeto tem como objetivo...
#   58 A primeira etapa...
#   59 Além disso, o projeto inclui a criação de um quiz...
#   60 Além disso, este projeto também terá uma seção de biblioteca virtual...
#   61 <new blank paragraph>
#   62 Premissas e Restrições (heading)

# 1d) Remove the "Diagrama" heading paragraph entirely.
$diagRng = $d.Content
$diagRng.Find.Execute("Diagrama", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$diagRng.Paragraphs.First.Range.Delete()

# 1e) Remove the "*****" run text (leaving the tab run on that line intact).
$d.Content.Find.Execute("*****", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ===========================================================================
# PHASE 2 - rewrite wording of the four surviving "Escopo" paragraphs.
# Re-fetched fresh by index (no structural edits happen from here on, so
# these indices are now stable).
# ===========================================================================

$old1 = "O meu projeto tem como objetivo criar um site que mostre a importância da leitura e como ela está presente na minha vida. Para alcançar esse objetivo, o escopo do projeto envolve diversas etapas, como coleta de dados dos usuários e análise desses dados para fornecer recomendações personalizadas de leitura."
$new1 = "Meu projeto tem como objetivo criar um site que demonstre a importância da leitura e sua presença significativa em minha vida. Para atingir essa meta, o escopo do projeto envolve diversas etapas, tais como coleta de dados dos usuários e análise dessas informações para fornecer recomendações de leitura personalizadas."
$d.Paragraphs.Item(57).Range.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "A primeira etapa é o levantamento de dados relevantes sobre livros e escritores famosos para que o site possa oferecer conteúdo de qualidade e útil para os usuários. Em seguida, será realizada a coleta de informações dos usuários, incluindo suas preferências de leitura e outras informações relevantes que possam ajudar a personalizar as recomendações."
$new2 = "A primeira etapa consiste em obter dados relevantes sobre livros e escritores renomados, a fim de oferecer conteúdo de qualidade e útil para os usuários. Em seguida, será realizada a coleta de informações dos usuários, incluindo suas preferências de leitura e outros dados pertinentes que possam contribuir para personalizar as recomendações."
$d.Paragraphs.Item(58).Range.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# This paragraph originally had multiple runs (incl. proofErr-wrapped text);
# replace the whole paragraph content (excluding the paragraph mark) with a
# single run carrying the new wording.
$p3 = $d.Paragraphs.Item(59)
$rng3 = $p3.Range
$rng3 = $d.Range($rng3.Start, $rng3.End - 1)
$rng3.Text = "Além disso, o projeto inclui uma comunidade em que os usuários podem compartilhar informações entre si e visualizar os dados de outros usuários."

$old4 = "Além disso, este projeto também terá uma seção de biblioteca virtual, onde o usuário poderá criar sua própria lista de livros que deseja ler ou já leu. Essa funcionalidade permitirá que o usuário mantenha um registro organizado de todos os livros que já leu ou ainda deseja ler."
$new4 = "Além disso, o projeto contará com uma seção de biblioteca virtual, em que o usuário poderá criar sua própria lista de livros desejados ou já lidos. Essa funcionalidade permitirá que o usuário mantenha um registro organizado de todas as obras que já foram lidas ou ainda desejam ser lidas."
$d.Paragraphs.Item(60).Range.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ===========================================================================
# PHASE 3 - apply the "Heading 1 style, rendered as normal Arial 11pt black
# text" direct formatting (pStyle Ttulo1 + run-level overrides) to the new
# blank paragraph and the four rewritten paragraphs.
# ===========================================================================

function Set-EscopoFormatting($para) {
    $para.Range.Style = $d.Styles.Item("Heading 1")
    $para.Range.Font.Name = "Arial"
    $para.Range.Font.NameFarEast = "Arial"
    $para.Range.Font.Bold = 0
    $para.Range.Font.BoldBi = 0
    $para.Range.Font.Size = 11
    $para.Range.Font.SizeBi = 11
    $para.Range.Font.Color = 0
    $para.Format.LineSpacingRule = 5
    $para.Format.LineSpacing = 18
}

Set-EscopoFormatting $d.Paragraphs.Item(56)

Set-EscopoFormatting $d.Paragraphs.Item(57)
$d.Paragraphs.Item(57).Format.FirstLineIndent = 35.4

Set-EscopoFormatting $d.Paragraphs.Item(58)
$d.Paragraphs.Item(58).Format.FirstLineIndent = 35.4

Set-EscopoFormatting $d.Paragraphs.Item(59)
$d.Paragraphs.Item(59).Format.FirstLineIndent = 35.4

Set-EscopoFormatting $d.Paragraphs.Item(60)
$d.Paragraphs.Item(60).Format.FirstLineIndent = 35.4
